# Auto-generated script to apply cryptos.xlsx diff changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "70.411.32"
Set-TextValue "E2" "  +0.89%  "
Set-TextValue "D3" "3.571.52"
Set-TextValue "E3" "  +0.07%  "
Set-TextValue "E4" "  -0.13%  "
Set-TextValue "E5" "  +2.96%  "
Set-TextValue "D6" "187.23"
Set-TextValue "E6" "  +0.45%  "
Set-TextValue "D7" "3.560.09"
Set-TextValue "E7" "  -0.19%  "
Set-TextValue "E8" "  +0.06%  "
Set-TextValue "E9" "  +0.02%  "
Set-TextValue "E10" "  +8.88%  "
Set-TextValue "D11" "0.649"
Set-TextValue "E11" "  -0.04%  "
Set-TextValue "D12" "54.81"
Set-TextValue "E12" "  -0.41%  "
Set-TextValue "E13" "  +1.61%  "
Set-TextValue "D14" "9.59"
Set-TextValue "E14" "  +0.55%  "
Set-TextValue "D15" "4.136.40"
Set-TextValue "E15" "  -0.33%  "
Set-TextValue "D16" "19.47"
Set-TextValue "E16" "  -0.55%  "
Set-TextValue "D17" "70.374.53"
Set-TextValue "E17" "  +0.78%  "
Set-TextValue "D18" "3.560.39"
Set-TextValue "E18" "  -0.45%  "
Set-TextValue "D19" "12.48"
Set-TextValue "E19" "  +0.02%  "
Set-TextValue "E20" "  -0.75%  "
Set-TextValue "D21" "552.29"
Set-TextValue "E21" "  +13.05%  "
Set-TextValue "E22" "  -0.57%  "
Set-TextValue "D23" "18.00"
Set-TextValue "E23" "  -8.22%  "
Set-TextValue "D24" "4.68"
Set-TextValue "E24" "  +8.72%  "
Set-TextValue "D25" "4.94"
Set-TextValue "E25" "  +0.07%  "
Set-TextValue "D26" "96.31"
Set-TextValue "E26" "  +0.58%  "
Set-TextValue "D27" "11.50"
Set-TextValue "E27" "  +4.42%  "
Set-TextValue "D28" "3.00"
Set-TextValue "E28" "  +1.90%  "
Set-TextValue "D29" "9.19"
Set-TextValue "E29" "  -0.70%  "
Set-TextValue "D30" "32.27"
Set-TextValue "E30" "  +1.78%  "
Set-TextValue "D31" "7.37"
Set-TextValue "E31" "  -1.57%  "
Set-TextValue "D32" "12.59"
Set-TextValue "E32" "  +4.88%  "
Set-TextValue "D33" "65.24"
Set-TextValue "E33" "  -2.20%  "
Set-TextValue "E34" "  -0.50%  "
Set-TextValue "B35" "Fetch.AI"
Set-TextValue "C35" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D35" "3.26"
Set-TextValue "E35" "  +6.42%  "
Set-TextValue "B36" "Bittensor"
Set-TextValue "C36" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D36" "553.44"
Set-TextValue "E36" "  -2.95%  "
Set-TextValue "E37" "  +7.30%  "
Set-TextValue "D38" "38.35"
Set-TextValue "E38" "  +0.31%  "
Set-TextValue "E39" "  +0.10%  "
Set-TextValue "E40" "  -3.09%  "
Set-TextValue "E41" "  -0.65%  "
Set-TextValue "D42" "3.382.16"
Set-TextValue "E42" "  +3.88%  "
Set-TextValue "B43" "dogwifhat"
Set-TextValue "C43" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D43" "3.13"
Set-TextValue "E43" "  -4.67%  "
Set-TextValue "B44" "Stacks"
Set-TextValue "C44" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D44" "3.39"
Set-TextValue "E44" "  -2.95%  "
Set-TextValue "E45" "  +3.25%  "
Set-TextValue "D47" "2.99"
Set-TextValue "E47" "  +0.09%  "
Set-TextValue "D48" "9.20"
Set-TextValue "E48" "  -4.36%  "
Set-TextValue "E49" "  +0.44%  "
Set-TextValue "D50" "1.00"
Set-TextValue "D51" "1.48"
Set-TextValue "E51" "  +22.71%  "

Write-Host "Applied all changes"
